$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 173
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 286
$ws.Range("K4").Value = 60
$ws.Range("L4").Value = 286
$ws.Range("M4").Value = 54
$ws.Range("N4").Value = -514

$ws.Range("H18").Value = 18845
$ws.Range("I18").Value = 18845
$ws.Range("K18").Value = 18845
$ws.Range("M18").Value = -18561

$ws.Range("H40").Value = 10106774
$ws.Range("J40").Value = 15879445
$ws.Range("L40").Value = 15879445
$ws.Range("N40").Value = -15879795

$ws.Range("H70").Value = 4780.905
$ws.Range("J70").Value = 6742.154
$ws.Range("L70").Value = 20226.462
$ws.Range("N70").Value = -20766.462

$ws.Range("H73").Value = 4780.905
$ws.Range("J73").Value = 6742.154
$ws.Range("L73").Value = 20226.462
$ws.Range("N73").Value = -22098.462

$ws.Range("H82").Value = 4203.2104
$ws.Range("I82").Value = 3415.4119
$ws.Range("J82").Value = 10899.5
$ws.Range("K82").Value = 10246.2357
$ws.Range("L82").Value = 32698.5
$ws.Range("M82").Value = -9840.235700000001
$ws.Range("N82").Value = -33510.5

$ws.Range("H85").Value = 4203.2104
$ws.Range("I85").Value = 3415.4119
$ws.Range("J85").Value = 10899.5
$ws.Range("K85").Value = 10246.2357
$ws.Range("L85").Value = 32698.5
$ws.Range("M85").Value = -8842.235700000001
$ws.Range("N85").Value = -35506.5

$ws.Range("H98").Value = 1202.2285
$ws.Range("I98").Value = 1202.2285
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1202.2285
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 295.7715000000001
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 2333.3333
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459

$ws.Range("H115").Value = 698.8889
$ws.Range("J115").Value = 1664
$ws.Range("L115").Value = 4992
$ws.Range("N115").Value = -8126

$ws.Range("H122").Value = 1202.2285
$ws.Range("I122").Value = 1202.2285
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3606.6855
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1156.6855
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1051871.5
$ws.Range("I2").Value = 1635445.6
$ws.Range("K2").Value = 1635445.6
$ws.Range("M2").Value = -1635332.6

$ws.Range("H61").Value = 55557052
$ws.Range("I61").Value = 55557052
$ws.Range("K61").Value = 55557052
$ws.Range("M61").Value = -55556840

$ws.Range("H116").Value = 1051871.5
$ws.Range("I116").Value = 1635445.6
$ws.Range("K116").Value = 1635445.6
$ws.Range("M116").Value = -1633151.6

$ws.Range("H136").Value = 55557052
$ws.Range("I136").Value = 55557052
$ws.Range("K136").Value = 166671156
$ws.Range("M136").Value = -166668606

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1051871.5
$ws.Range("I3").Value = 1635445.6
$ws.Range("K3").Value = 1635445.6
$ws.Range("M3").Value = -1635331.6

$ws.Range("H20").Value = 1641.1
$ws.Range("I20").Value = 1592.7
$ws.Range("K20").Value = 1592.7
$ws.Range("M20").Value = -1345.7

$ws.Range("H107").Value = 64240.375
$ws.Range("I107").Value = 1687.6364
$ws.Range("J107").Value = 201856.4
$ws.Range("K107").Value = 1687.6364
$ws.Range("L107").Value = 201856.4
$ws.Range("M107").Value = 232.3635999999999
$ws.Range("N107").Value = -205696.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 5367
$ws.Range("I39").Value = 5367
$ws.Range("K39").Value = 5367
$ws.Range("M39").Value = -4976

$ws.Range("H49").Value = 5367
$ws.Range("I49").Value = 5367
$ws.Range("K49").Value = 5367
$ws.Range("M49").Value = -5185

$ws.Range("H132").Value = 29413462
$ws.Range("I132").Value = 33334964
$ws.Range("K132").Value = 100004892
$ws.Range("M132").Value = -100002362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31.31579
$ws.Range("I2").Value = 21.214285
$ws.Range("K2").Value = 127.28571
$ws.Range("M2").Value = -14.28570999999999

$ws.Range("H3").Value = 10376.625
$ws.Range("I3").Value = 9716.143
$ws.Range("K3").Value = 29148.429
$ws.Range("M3").Value = -29036.429

$ws.Range("H4").Value = 4601608
$ws.Range("I4").Value = 6053537
$ws.Range("K4").Value = 18160611
$ws.Range("M4").Value = -18160499

$ws.Range("H38").Value = 70.875
$ws.Range("I38").Value = 55.8
$ws.Range("J38").Value = 96
$ws.Range("K38").Value = 167.4
$ws.Range("L38").Value = 288
$ws.Range("M38").Value = 179.6
$ws.Range("N38").Value = -982

$ws.Range("H121").Value = 150323
$ws.Range("I121").Value = 200117.8
$ws.Range("J121").Value = 67331.664
$ws.Range("K121").Value = 600353.3999999999
$ws.Range("L121").Value = 201994.992
$ws.Range("M121").Value = -599043.3999999999
$ws.Range("N121").Value = -204614.992

$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -15880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3194.1738
$ws.Range("I80").Value = 3310.8
$ws.Range("J80").Value = 2975.5
$ws.Range("K80").Value = 3310.8
$ws.Range("L80").Value = 2975.5
$ws.Range("M80").Value = -2312.8
$ws.Range("N80").Value = -4971.5

$ws.Range("H83").Value = 3194.1738
$ws.Range("I83").Value = 3310.8
$ws.Range("J83").Value = 2975.5
$ws.Range("K83").Value = 16554
$ws.Range("L83").Value = 14877.5
$ws.Range("M83").Value = -11562
$ws.Range("N83").Value = -24861.5

$ws.Range("H113").Value = 32029
$ws.Range("I113").Value = 35567.43
$ws.Range("K113").Value = 35567.43
$ws.Range("M113").Value = -33397.43

$ws.Range("H122").Value = 5460.7188
$ws.Range("I122").Value = 4379.7617
$ws.Range("K122").Value = 13139.2851
$ws.Range("M122").Value = -10689.2851

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1733
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 1500
$ws.Range("M82").Value = -1139

$ws.Range("H85").Value = 1733
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 1500
$ws.Range("M85").Value = -252

$ws.Range("H93").Value = 2575.7144
$ws.Range("I93").Value = 1527
$ws.Range("J93").Value = 3362.25
$ws.Range("K93").Value = 1527
$ws.Range("L93").Value = 3362.25
$ws.Range("M93").Value = -279
$ws.Range("N93").Value = -5858.25

$ws.Range("H136").Value = 2117.8914
$ws.Range("I136").Value = 1743.4
$ws.Range("J136").Value = 2299.0967
$ws.Range("K136").Value = 5230.200000000001
$ws.Range("L136").Value = 6897.2901
$ws.Range("M136").Value = -2680.200000000001
$ws.Range("N136").Value = -11997.2901

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1004.46155
$ws.Range("I126").Value = 926.8
$ws.Range("K126").Value = 2780.4
$ws.Range("M126").Value = -310.3999999999996
